# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.651.59"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.843.67"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  -1.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.37"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4304"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3732"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07334"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8746"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.51"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "1.861.11"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.707"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.434"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07126"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.95"
$ws.Range("E16").Value = "  +4.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.017"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008963"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.013"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.42"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").Value = "27.664.32"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.213"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.06"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D24").Value = "2.084.77"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.011"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.78"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.58"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.138"
$ws.Range("E28").Value = "  +7.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.353"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.86"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08951"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.224"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7729"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.538"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -3.61%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.135"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05316"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01969"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.243"
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.886"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5106"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1681"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.742"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.71"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "109.38"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4724"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06491"
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.687"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.835"
$ws.Range("E51").Value = "  -4.50%  "
